# Case_4_42 vm_pu.xlsx results update ("case with 380 kV done")
# Re-run of the power flow case: slack/reference voltage setpoint on
# column B changed from 1.05 pu to 1.02 pu, and the resulting bus voltage
# magnitudes (columns C-F and I-N, rows 2-25) are updated to the newly
# computed values. Column G (=1) is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("B3").Value = 1.02
$ws.Range("B4").Value = 1.02
$ws.Range("B5").Value = 1.02
$ws.Range("B6").Value = 1.02
$ws.Range("B7").Value = 1.02
$ws.Range("B8").Value = 1.02
$ws.Range("B9").Value = 1.02
$ws.Range("B10").Value = 1.02
$ws.Range("B11").Value = 1.02
$ws.Range("B12").Value = 1.02
$ws.Range("B13").Value = 1.02
$ws.Range("B14").Value = 1.02
$ws.Range("B15").Value = 1.02
$ws.Range("B16").Value = 1.02
$ws.Range("B17").Value = 1.02
$ws.Range("B18").Value = 1.02
$ws.Range("B19").Value = 1.02
$ws.Range("B20").Value = 1.02
$ws.Range("B21").Value = 1.02
$ws.Range("B22").Value = 1.02
$ws.Range("B23").Value = 1.02
$ws.Range("B24").Value = 1.02
$ws.Range("B25").Value = 1.02
$ws.Range("C2").Value = 1.050322759995683
$ws.Range("C3").Value = 1.05176301652193
$ws.Range("C4").Value = 1.052693844550915
$ws.Range("C5").Value = 1.053084904096581
$ws.Range("C6").Value = 1.053150549615361
$ws.Range("C7").Value = 1.052699070923833
$ws.Range("C8").Value = 1.050809735882487
$ws.Range("C9").Value = 1.047471674065325
$ws.Range("C10").Value = 1.045239983890099
$ws.Range("C11").Value = 1.044272051897722
$ws.Range("C12").Value = 1.043912271545931
$ws.Range("C13").Value = 1.043989456966979
$ws.Range("C14").Value = 1.044242317430961
$ws.Range("C15").Value = 1.044398080108169
$ws.Range("C16").Value = 1.045304188085801
$ws.Range("C17").Value = 1.045872133178482
$ws.Range("C18").Value = 1.046203252525181
$ws.Range("C19").Value = 1.046316129872769
$ws.Range("C20").Value = 1.04581121396894
$ws.Range("C21").Value = 1.044167863215177
$ws.Range("C22").Value = 1.04313318868702
$ws.Range("C23").Value = 1.043681827778527
$ws.Range("C24").Value = 1.045838741211841
$ws.Range("C25").Value = 1.048335727432832
$ws.Range("D2").Value = 1.049411262395481
$ws.Range("D3").Value = 1.050497871104509
$ws.Range("D4").Value = 1.051199771170514
$ws.Range("D5").Value = 1.051494564410726
$ws.Range("D6").Value = 1.051544044856895
$ws.Range("D7").Value = 1.051203711329584
$ws.Range("D8").Value = 1.049778739374728
$ws.Range("D9").Value = 1.0472583325996
$ws.Range("D10").Value = 1.045571486368116
$ws.Range("D11").Value = 1.044839449921481
$ws.Range("D12").Value = 1.044567290731848
$ws.Range("D13").Value = 1.044625681112624
$ws.Range("D14").Value = 1.044816958249923
$ws.Range("D15").Value = 1.044934777368
$ws.Range("D16").Value = 1.045620034662176
$ws.Range("D17").Value = 1.04604944116474
$ws.Range("D18").Value = 1.04629975064242
$ws.Range("D19").Value = 1.046385073375948
$ws.Range("D20").Value = 1.046003386088382
$ws.Range("D21").Value = 1.044760638791857
$ws.Range("D22").Value = 1.043977835515686
$ws.Range("D23").Value = 1.044392952308979
$ws.Range("D24").Value = 1.046024196877535
$ws.Range("D25").Value = 1.047911059444123
$ws.Range("E2").Value = 1.063582166126474
$ws.Range("E3").Value = 1.065001567912547
$ws.Range("E4").Value = 1.065919256169349
$ws.Range("E5").Value = 1.066304875915054
$ws.Range("E6").Value = 1.066369612981181
$ws.Range("E7").Value = 1.065924409525621
$ws.Range("E8").Value = 1.064062020179339
$ws.Range("E9").Value = 1.060774189349334
$ws.Range("E10").Value = 1.058577880917069
$ws.Range("E11").Value = 1.057625729078412
$ws.Range("E12").Value = 1.05727188010039
$ws.Range("E13").Value = 1.057347790043525
$ws.Range("E14").Value = 1.057596483469432
$ws.Range("E15").Value = 1.057749687988349
$ws.Range("E16").Value = 1.05864104759377
$ws.Range("E17").Value = 1.059199865088992
$ws.Range("E18").Value = 1.059525704750975
$ws.Range("E19").Value = 1.059636789348352
$ws.Range("E20").Value = 1.059139920597034
$ws.Range("E21").Value = 1.057523254405315
$ws.Range("E22").Value = 1.056505762500771
$ws.Range("E23").Value = 1.057045254080017
$ws.Range("E24").Value = 1.059167007270205
$ws.Range("E25").Value = 1.061624926279094
$ws.Range("F2").Value = 1.070183829371126
$ws.Range("F3").Value = 1.071651138091094
$ws.Range("F4").Value = 1.072599851997329
$ws.Range("F5").Value = 1.072998521525962
$ws.Range("F6").Value = 1.073065450101396
$ws.Range("F7").Value = 1.072605179699235
$ws.Range("F8").Value = 1.070679868407583
$ws.Range("F9").Value = 1.067281355515391
$ws.Range("F10").Value = 1.065011376165843
$ws.Range("F11").Value = 1.06402734876061
$ws.Range("F12").Value = 1.063661663186908
$ws.Range("F13").Value = 1.063740111969935
$ws.Range("F14").Value = 1.06399712466691
$ws.Range("F15").Value = 1.064155455407204
$ws.Range("F16").Value = 1.065076658811328
$ws.Range("F17").Value = 1.065654202889008
$ws.Range("F18").Value = 1.065990967824232
$ws.Range("F19").Value = 1.066105778068075
$ws.Range("F20").Value = 1.065592248968263
$ws.Range("F21").Value = 1.063921445692812
$ws.Range("F22").Value = 1.062869935248321
$ws.Range("F23").Value = 1.063427458929436
$ws.Range("F24").Value = 1.065620243608477
$ws.Range("F25").Value = 1.068160684766334
$ws.Range("I2").Value = 1.03941775864605
$ws.Range("I3").Value = 1.039739019837312
$ws.Range("I4").Value = 1.039945024042573
$ws.Range("I5").Value = 1.040031181231499
$ws.Range("I6").Value = 1.040045621234542
$ws.Range("I7").Value = 1.039946177033638
$ws.Range("I8").Value = 1.039526719139045
$ws.Range("I9").Value = 1.038773166764646
$ws.Range("I10").Value = 1.038261012456692
$ws.Range("I11").Value = 1.03803690134257
$ws.Range("I12").Value = 1.037953302281098
$ws.Range("I13").Value = 1.037971250623497
$ws.Range("I14").Value = 1.038029998255098
$ws.Range("I15").Value = 1.03806614761528
$ws.Range("I16").Value = 1.038275836405164
$ws.Range("I17").Value = 1.038406739492509
$ws.Range("I18").Value = 1.03848286692864
$ws.Range("I19").Value = 1.038508786120806
$ws.Range("I20").Value = 1.03839271823268
$ws.Range("I21").Value = 1.038012708335383
$ws.Range("I22").Value = 1.037771730690113
$ws.Range("I23").Value = 1.037899672491702
$ws.Range("I24").Value = 1.038399054535277
$ws.Range("I25").Value = 1.038969695970889
$ws.Range("J2").Value = 1.055356524378297
$ws.Range("J3").Value = 1.056444019774661
$ws.Range("J4").Value = 1.057146162242889
$ws.Range("J5").Value = 1.057440978775224
$ws.Range("J6").Value = 1.057490458542933
$ws.Range("J7").Value = 1.057150103021205
$ws.Range("J8").Value = 1.055724370052312
$ws.Range("J9").Value = 1.053200052899261
$ws.Range("J10").Value = 1.051508828420443
$ws.Range("J11").Value = 1.050774466397818
$ws.Range("J12").Value = 1.050501377999173
$ws.Range("J13").Value = 1.050559970670487
$ws.Range("J14").Value = 1.050751899256163
$ws.Range("J15").Value = 1.050870111064197
$ws.Range("J16").Value = 1.051557521930782
$ws.Range("J17").Value = 1.051988163674292
$ws.Range("J18").Value = 1.052239152236214
$ws.Range("J19").Value = 1.052324699520705
$ws.Range("J20").Value = 1.051941980357514
$ws.Range("J21").Value = 1.050695389774364
$ws.Range("J22").Value = 1.049909791357848
$ws.Range("J23").Value = 1.050326426079573
$ws.Range("J24").Value = 1.051962849220903
$ws.Range("J25").Value = 1.053854099662696
$ws.Range("K2").Value = 1.052168395498476
$ws.Range("K3").Value = 1.053066536422366
$ws.Range("K4").Value = 1.053645895310644
$ws.Range("K5").Value = 1.053889030822444
$ws.Range("K6").Value = 1.053929829410448
$ws.Range("K7").Value = 1.053649145772532
$ws.Range("K8").Value = 1.052472300939046
$ws.Range("K9").Value = 1.050384616394462
$ws.Range("K10").Value = 1.048983232135064
$ws.Range("K11").Value = 1.048374090544281
$ws.Range("K12").Value = 1.048147473247338
$ws.Range("K13").Value = 1.048196099563173
$ws.Range("K14").Value = 1.04835536555555
$ws.Range("K15").Value = 1.048453447415707
$ws.Range("K16").Value = 1.049023609321273
$ws.Range("K17").Value = 1.049380629064092
$ws.Range("K18").Value = 1.049588647848624
$ws.Range("K19").Value = 1.049659538884947
$ws.Range("K20").Value = 1.049342347508283
$ws.Range("K21").Value = 1.048308475526011
$ws.Range("K22").Value = 1.047656383077555
$ws.Range("K23").Value = 1.048002266023287
$ws.Range("K24").Value = 1.049359645990387
$ws.Range("K25").Value = 1.050926007298675
$ws.Range("L2").Value = 1.066300415962667
$ws.Range("L3").Value = 1.06753334085137
$ws.Range("L4").Value = 1.068329840609129
$ws.Range("L5").Value = 1.06866438639577
$ws.Range("L6").Value = 1.068720540497408
$ws.Range("L7").Value = 1.068334312009633
$ws.Range("L8").Value = 1.066717357768817
$ws.Range("L9").Value = 1.063858013445727
$ws.Range("L10").Value = 1.061944704540719
$ws.Range("L11").Value = 1.061114471332078
$ws.Range("L12").Value = 1.060805815524362
$ws.Range("L13").Value = 1.060872035570203
$ws.Range("L14").Value = 1.061088963294929
$ws.Range("L15").Value = 1.061222583633599
$ws.Range("L16").Value = 1.061999766826458
$ws.Range("L17").Value = 1.062486798372541
$ws.Range("L18").Value = 1.062770706445907
$ws.Range("L19").Value = 1.062867483249567
$ws.Range("L20").Value = 1.062434562041261
$ws.Range("L21").Value = 1.061025090995868
$ws.Range("L22").Value = 1.060137334083248
$ws.Range("L23").Value = 1.060608101366344
$ws.Range("L24").Value = 1.062458165913957
$ws.Range("L25").Value = 1.06459844723202
$ws.Range("M2").Value = 1.072884345364755
$ws.Range("M3").Value = 1.074166356283882
$ws.Range("M4").Value = 1.074994651775184
$ws.Range("M5").Value = 1.07534257275961
$ws.Range("M6").Value = 1.075400973106131
$ws.Range("M7").Value = 1.074999301863008
$ws.Range("M8").Value = 1.073317869138118
$ws.Range("M9").Value = 1.070345156113667
$ws.Range("M10").Value = 1.068356425141417
$ws.Range("M11").Value = 1.067493568045802
$ws.Range("M12").Value = 1.067172799216022
$ws.Range("M13").Value = 1.067241617326134
$ws.Range("M14").Value = 1.067467058642935
$ws.Range("M15").Value = 1.067605925134917
$ws.Range("M16").Value = 1.068413653272309
$ws.Range("M17").Value = 1.068919853900213
$ws.Range("M18").Value = 1.069214946305501
$ws.Range("M19").Value = 1.069315537248938
$ws.Range("M20").Value = 1.068865560564951
$ws.Range("M21").Value = 1.067400679168557
$ws.Range("M22").Value = 1.066478112052382
$ws.Range("M23").Value = 1.066967330257182
$ws.Range("M24").Value = 1.068890093898921
$ws.Range("M25").Value = 1.071114869900749
$ws.Range("N2").Value = 1.022102396198667
$ws.Range("N3").Value = 1.022478027086057
$ws.Range("N4").Value = 1.022720220084916
$ws.Range("N5").Value = 1.022821832257814
$ws.Range("N6").Value = 1.022838881362934
$ws.Range("N7").Value = 1.022721578638533
$ws.Range("N8").Value = 1.022229522422784
$ws.Range("N9").Value = 1.021355766580713
$ws.Range("N10").Value = 1.020768675971242
$ws.Range("N11").Value = 1.020513351260738
$ws.Range("N12").Value = 1.020418343664389
$ws.Range("N13").Value = 1.020438730775861
$ws.Range("N14").Value = 1.020505501353517
$ws.Range("N15").Value = 1.020546618511586
$ws.Range("N16").Value = 1.020785597481092
$ws.Range("N17").Value = 1.020935204086383
$ws.Range("N18").Value = 1.021022360052063
$ws.Range("N19").Value = 1.021052059882227
$ws.Range("N20").Value = 1.020919163797635
$ws.Range("N21").Value = 1.020485843751967
$ws.Range("N22").Value = 1.020212421706706
$ws.Range("N23").Value = 1.020357461066543
$ws.Range("N24").Value = 1.020926412044527
$ws.Range("N25").Value = 1.021582455295699
